$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.877.99"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.821.37"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "'240.73"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").Value = "'0.6149"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").Value = "'0.9957"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.07388"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "'0.2917"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'22.89"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "'0.07616"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "1.820.18"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "'4.964"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "'0.6683"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'82.39"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "'0.000009016"
$ws.Range("E16").Value = "  -6.54%  "
$ws.Range("D17").Value = "'5.852"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "28.892.63"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "2.115.62"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "'241.16"
$ws.Range("E20").Value = "  +7.97%  "
$ws.Range("D21").Value = "'12.62"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "'0.9957"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'7.176"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").Value = "'157.69"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("D26").Value = "'0.1409"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'8.431"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'17.75"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "'1.482"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("D30").Value = "'0.05563"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.083"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.099"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "'1.205"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("D34").Value = "'1.825"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "'0.7367"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "'2.622"
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("D38").Value = "'2.752"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "1.201.34"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "'6.345"
$ws.Range("E41").Value = "  -4.71%  "
$ws.Range("D42").Value = "'0.8926"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'0.9939"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "2.004.76"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").Value = "'100.93"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "'64.91"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5059"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").Value = "'0.4032"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000117"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").Value = "'9.061"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "'0.05794"
$ws.Range("E51").Value = "  +0.26%  "
